$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the SqlIP value (E2) to the new address
$ws.Range("E2").Value = "192.168.0.24"

# Update the active cell selection from G6 to H6
$ws.Range("H6").Select()
